# "Added more Redfin TestCases with Assertions"
# Append 11 more city names to the Redfin test-data sheet (column A),
# continuing directly below the existing list (rows 2-31 -> now 2-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Redfin")

$cities = @(
    "Destin",
    "Fort Lauderdale",
    "Jacksonville",
    "Fort Myers",
    "Charleston",
    "Myrtle Beach",
    "Asheville",
    "Durham",
    "Raleigh",
    "Charlotte",
    "Louisville"
)

$startRow = 32
for ($i = 0; $i -lt $cities.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $cities[$i]
}

# Move/select the last appended cell, matching the saved selection state.
$lastRow = $startRow + $cities.Length - 1
$ws.Range("A$lastRow").Select()

# Scroll the window so row 11 is at the top of the visible area (best effort;
# matches the author's saved scroll position).
$excel.ActiveWindow.ScrollRow = 11
